{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change summary (matches the target OOXML diff):\n//  1. Strip the \"Heading2\" paragraph style from the five section headings\n//     (\"Introduction\", \"Cognitive Psychology Perspective\",\n//     \"Social Psychology Perspective\", \"Comparative Analysis of\n//     Perspectives\", \"Conclusion\") so they fall back to the default\n//     (Normal) paragraph style.\n//  2. Swap placeholder/author citations for the final reference keys in\n//     five body paragraphs (citation-check pass).\n\nconst headingTexts = [\n  \"Introduction\",\n  \"Cognitive Psychology Perspective\",\n  \"Social Psychology Perspective\",\n  \"Comparative Analysis of Perspectives\",\n  \"Conclusion\",\n];\n\nconst bodyReplacements = [\n  { prefix: \"Cognitive Psychology is a field dedicated to understanding mental processes\", text: \"Cognitive Psychology is a field dedicated to understanding mental processes such as perception, memory, and reasoning, and it extensively investigates how these processes are influenced by social media. Psychologists in this domain explore the effects of social media on cognitive functions, particularly focusing on how digital interactions can foster aggressive behavior and bullying. One pivotal study in \\\"A Cognitive Psychology of Mass Communication\\\" highlights how exposure to aggressive content on social media can distort perceptions of reality and increase hostile cognitive biases (Lee 208). This perspective posits that social media platforms, by manipulating the accessibility of aggressive stimuli, can alter cognitive processing, resulting in desensitization to violence. Moreover, textbooks and recent journal articles consistently support the notion that these cognitive distortions are exacerbated by the pervasive and interactive nature of social media, which amplifies the effects on users' thought patterns and subsequent behaviors (Lee 208).\" },\n  { prefix: \"Furthermore, cognitive psychologists examining social media's impact\", text: \"Furthermore, cognitive psychologists examining social media's impact often pose questions about the specific cognitive biases that are exacerbated by virtual interactions. They might explore how social media platforms influence users\u2019 attention and memory, particularly regarding the spread of misinformation and disinformation (Pearse et al., 2001). Researchers in this field are likely to conduct experiments that test the effects of repeated exposure to aggressive content on cognitive processing, such as the accessibility of violent schemas and its impact on real-world behavior. Additionally, these psychologists investigate the phenomenon of social media rumination, where users continuously think about online interactions, potentially disrupting positive cognitive distractions and leading to increased distress and aggression (Pearse et al., 2001). By addressing these questions through empirical studies, cognitive psychologists aim to uncover the mechanisms by which digital environments shape cognitive functions, ultimately informing strategies to mitigate negative psychological outcomes.\" },\n  { prefix: \"Social Psychology investigates how social media influences behavior\", text: \"Social Psychology investigates how social media influences behavior by examining group dynamics and interpersonal interactions, particularly focusing on phenomena such as aggressive behavior and bullying. Social psychologists interpret the effects of social media through the lens of social norms, peer influence, and the role of anonymity in exacerbating aggressive tendencies. A key study highlights how social media platforms can intensify aggressive behavior due to reduced social cues and anonymity, which often lead to cyberbullying among adolescents (Ref-u156544). By integrating socio-psychological frameworks, researchers can explore how these platforms amplify aggression, as well as the potential for empathy and moral disengagement to moderate these effects. Furthermore, literature such as \\\"The Social Psychology of Aggression\\\" underscores the significance of situational factors and media exposure in fostering aggression, providing insights into strategies for prevention and intervention (Ref-u156544).\" },\n  { prefix: \"Moreover, social psychologists are particularly interested in exploring\", text: \"Moreover, social psychologists are particularly interested in exploring how peer influence and conformity impact behavior on social media platforms. They might investigate how social norms are established and perpetuated in online communities, focusing on the role of anonymity and reduced social cues in fostering aggressive interactions, such as cyberbullying (Ref-s948613). Researchers in this field often employ both qualitative and quantitative methods to examine the dynamics of peer pressure and conformity, seeking to understand how these social forces contribute to bullying and aggression across various digital contexts (Ref-s948613). This includes analyzing cross-national data to assess cultural variations in social media use and its relationship with aggressive behaviors, thereby providing insights into different socio-cultural contexts (Ref-s948613). By addressing these questions, social psychologists aim to develop targeted intervention strategies that can mitigate the negative impacts of social media and promote healthier virtual interactions.\" },\n  { prefix: \"In examining the psychological impacts of social media, both Cognitive\", text: \"In examining the psychological impacts of social media, both Cognitive and Social Psychology offer unique yet complementary insights into phenomena such as aggression and bullying. Cognitive Psychology focuses on how mental processes are altered by digital environments, emphasizing the role of cognitive biases and memory distortions (Ref-u482754). Conversely, Social Psychology prioritizes the influence of social norms, peer pressure, and the anonymity facilitated by online interactions, which can exacerbate aggressive behaviors like cyberbullying (Ref-u482754). Despite these differences, both perspectives acknowledge the inherent complexity of online interactions and recognize that social media can act as a catalyst for both victimization and perpetration of bullying, particularly when problematic use is involved (Ref-u482754). By integrating cognitive insights with social frameworks, researchers can develop more comprehensive strategies to address the multifaceted nature of social media's impact, thereby fostering healthier digital environments and reducing the prevalence of aggression and bullying.\" },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n\n  if (headingTexts.indexOf(text) !== -1) {\n    // Remove the Heading2 style -> paragraph reverts to the default style.\n    para.style = \"Normal\";\n    continue;\n  }\n\n  for (let r = 0; r < bodyReplacements.length; r++) {\n    const rep = bodyReplacements[r];\n    if (para.text.indexOf(rep.prefix) === 0) {\n      para.insertText(rep.text, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Change summary (matches the target OOXML diff):\n#  1. Strip the \"Heading2\" paragraph style from the five section headings\n#     (\"Introduction\", \"Cognitive Psychology Perspective\",\n#     \"Social Psychology Perspective\", \"Comparative Analysis of\n#     Perspectives\", \"Conclusion\") so they fall back to the default\n#     (Normal) paragraph style.\n#  2. Swap placeholder/author citations for the final reference keys in\n#     five body paragraphs (citation-check pass).\n\n$d = $word.ActiveDocument\n\n$headingTexts = @(\n    \"Introduction\",\n    \"Cognitive Psychology Perspective\",\n    \"Social Psychology Perspective\",\n    \"Comparative Analysis of Perspectives\",\n    \"Conclusion\"\n)\n\n$bodyReplacements = @(\n    @{ Prefix = 'Cognitive Psychology is a field dedicated to understanding mental processes'; Text = 'Cognitive Psychology is a field dedicated to understanding mental processes such as perception, memory, and reasoning, and it extensively investigates how these processes are influenced by social media. Psychologists in this domain explore the effects of social media on cognitive functions, particularly focusing on how digital interactions can foster aggressive behavior and bullying. One pivotal study in \"A Cognitive Psychology of Mass Communication\" highlights how exposure to aggressive content on social media can distort perceptions of reality and increase hostile cognitive biases (Lee 208). This perspective posits that social media platforms, by manipulating the accessibility of aggressive stimuli, can alter cognitive processing, resulting in desensitization to violence. Moreover, textbooks and recent journal articles consistently support the notion that these cognitive distortions are exacerbated by the pervasive and interactive nature of social media, which amplifies the effects on users'' thought patterns and subsequent behaviors (Lee 208).' },\n    @{ Prefix = 'Furthermore, cognitive psychologists examining social media''s impact'; Text = 'Furthermore, cognitive psychologists examining social media''s impact often pose questions about the specific cognitive biases that are exacerbated by virtual interactions. They might explore how social media platforms influence users\u2019 attention and memory, particularly regarding the spread of misinformation and disinformation (Pearse et al., 2001). Researchers in this field are likely to conduct experiments that test the effects of repeated exposure to aggressive content on cognitive processing, such as the accessibility of violent schemas and its impact on real-world behavior. Additionally, these psychologists investigate the phenomenon of social media rumination, where users continuously think about online interactions, potentially disrupting positive cognitive distractions and leading to increased distress and aggression (Pearse et al., 2001). By addressing these questions through empirical studies, cognitive psychologists aim to uncover the mechanisms by which digital environments shape cognitive functions, ultimately informing strategies to mitigate negative psychological outcomes.' },\n    @{ Prefix = 'Social Psychology investigates how social media influences behavior'; Text = 'Social Psychology investigates how social media influences behavior by examining group dynamics and interpersonal interactions, particularly focusing on phenomena such as aggressive behavior and bullying. Social psychologists interpret the effects of social media through the lens of social norms, peer influence, and the role of anonymity in exacerbating aggressive tendencies. A key study highlights how social media platforms can intensify aggressive behavior due to reduced social cues and anonymity, which often lead to cyberbullying among adolescents (Ref-u156544). By integrating socio-psychological frameworks, researchers can explore how these platforms amplify aggression, as well as the potential for empathy and moral disengagement to moderate these effects. Furthermore, literature such as \"The Social Psychology of Aggression\" underscores the significance of situational factors and media exposure in fostering aggression, providing insights into strategies for prevention and intervention (Ref-u156544).' },\n    @{ Prefix = 'Moreover, social psychologists are particularly interested in exploring'; Text = 'Moreover, social psychologists are particularly interested in exploring how peer influence and conformity impact behavior on social media platforms. They might investigate how social norms are established and perpetuated in online communities, focusing on the role of anonymity and reduced social cues in fostering aggressive interactions, such as cyberbullying (Ref-s948613). Researchers in this field often employ both qualitative and quantitative methods to examine the dynamics of peer pressure and conformity, seeking to understand how these social forces contribute to bullying and aggression across various digital contexts (Ref-s948613). This includes analyzing cross-national data to assess cultural variations in social media use and its relationship with aggressive behaviors, thereby providing insights into different socio-cultural contexts (Ref-s948613). By addressing these questions, social psychologists aim to develop targeted intervention strategies that can mitigate the negative impacts of social media and promote healthier virtual interactions.' },\n    @{ Prefix = 'In examining the psychological impacts of social media, both Cognitive'; Text = 'In examining the psychological impacts of social media, both Cognitive and Social Psychology offer unique yet complementary insights into phenomena such as aggression and bullying. Cognitive Psychology focuses on how mental processes are altered by digital environments, emphasizing the role of cognitive biases and memory distortions (Ref-u482754). Conversely, Social Psychology prioritizes the influence of social norms, peer pressure, and the anonymity facilitated by online interactions, which can exacerbate aggressive behaviors like cyberbullying (Ref-u482754). Despite these differences, both perspectives acknowledge the inherent complexity of online interactions and recognize that social media can act as a catalyst for both victimization and perpetration of bullying, particularly when problematic use is involved (Ref-u482754). By integrating cognitive insights with social frameworks, researchers can develop more comprehensive strategies to address the multifaceted nature of social media''s impact, thereby fostering healthier digital environments and reducing the prevalence of aggression and bullying.' }\n)\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    $raw = $p.Range.Text\n    $trimmed = $raw.TrimEnd([char]13, [char]7)\n\n    if ($headingTexts -contains $trimmed) {\n        # Remove the Heading2 style -> paragraph reverts to the default style.\n        $p.Style = \"Normal\"\n        continue\n    }\n\n    foreach ($rep in $bodyReplacements) {\n        if ($trimmed.StartsWith($rep.Prefix)) {\n            $p.Range.Text = $rep.Text\n            break\n        }\n    }\n}\n"}
